# Apply MaxInvest / ExisUnits updates on the "Power Storage" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power Storage")

# ExisUnits (column E) changes
$ws.Range("E7").Value = 7
$ws.Range("E9").Value = 0
$ws.Range("E10").Value = 33

# MaxInvest (column S) changes
$ws.Range("S7").Value = 15
$ws.Range("S8").Value = 15
$ws.Range("S9").Value = 15
$ws.Range("S10").Value = 15
$ws.Range("S11").Value = 15
